$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Valor Mora" amount
$ws.Range("E11").Value = 120134

# 2. Update "Cant. Periodos" count
$ws.Range("F13").Value = 2

# 3. Insert a new data row (2509 / 68000) below the existing data row (16),
#    copying its formatting, then push the two footer rows down by one.
$ws.Rows("16:16").Copy()
$ws.Rows("17:17").Insert()

$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 68000
